# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, as captured in the "Updated cryptos list ... with GitHub
# Actions" commit. Price cells in column D are stored as plain text
# (European-style "thousand.decimal" strings like "38.771.68" aren't valid
# numbers, and short ones like "227.40"/"0.120" must keep their trailing
# zero), so every D write is preceded by forcing the cell to Text format -
# otherwise Excel's COM layer would silently coerce them to numbers and
# drop significant trailing/duplicate-dot formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.771.68"
$ws.Range("E2").Value = "  -0.09%  "
# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.33"
$ws.Range("E3").Value = "  +0.24%  "
# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "
# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.40"
$ws.Range("E5").Value = "  -0.45%  "
# Row 6 - XRP
$ws.Range("E6").Value = "  +0.39%  "
# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.77"
$ws.Range("E7").Value = "  +2.26%  "
# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "
# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.97%  "
# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").Value = "  +0.66%  "
# Row 11 - TRON
$ws.Range("E11").Value = "  -1.00%  "
# Row 12 - Chainlink
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.76"
$ws.Range("E12").Value = "  +5.33%  "
# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.413.25"
$ws.Range("E13").Value = "  +0.22%  "
# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.98"
$ws.Range("E14").Value = "  -0.68%  "
# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.35%  "
# Row 16 - Polkadot
$ws.Range("E16").Value = "  +0.09%  "
# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.058.16"
$ws.Range("E17").Value = "  -1.98%  "
# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.901.89"
# Row 19 - Litecoin
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.59"
$ws.Range("E19").Value = "  -0.64%  "
# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  -0.71%  "
# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +0.99%  "
# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.86"
$ws.Range("E22").Value = "  +0.21%  "
# Row 23 - Dai
$ws.Range("E23").Value = "  +0.04%  "
# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -2.66%  "
# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.84%  "
# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.63"
$ws.Range("E26").Value = "  +1.30%  "
# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.15"
$ws.Range("E27").Value = "  -0.23%  "
# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.38%  "
# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +0.58%  "
# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.34"
$ws.Range("E30").Value = "  +0.78%  "
# Row 31 - WEMIXToken
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.56"
$ws.Range("E31").Value = "  +9.41%  "
# Row 32 - Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  +0.13%  "
# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.55"
$ws.Range("E33").Value = "  +1.45%  "

# Row 34/35 swap: row 34 was InternetComputer(DFINITY), now THORChain;
# row 35 was THORChain, now InternetComputer(DFINITY).
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.15"
$ws.Range("E34").Value = "  +12.04%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.76"
$ws.Range("E35").Value = "  +0.28%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0613"
$ws.Range("E36").Value = "  -0.06%  "
# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -1.50%  "
# Row 38 - RenderToken
$ws.Range("E38").Value = "  -0.33%  "
# Row 39 - BinanceUSD
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40/41 swap: row 40 was InjectiveProtocol, now VeChain;
# row 41 was VeChain, now InjectiveProtocol.
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0229"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.95"
$ws.Range("E41").Value = "  -2.49%  "

# Row 42 - Aave
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.84"
$ws.Range("E42").Value = "  +0.19%  "
# Row 43 - Maker
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.523.95"
$ws.Range("E43").Value = "  -1.17%  "
# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +7.10%  "
# Row 45 - HuobiToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("E45").Value = "  -0.59%  "
# Row 46 - FraxShare
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.76"
$ws.Range("E46").Value = "  +1.06%  "
# Row 47 - Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0910"
$ws.Range("E47").Value = "  -1.76%  "
# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +4.67%  "
# Row 49 - FTXToken
$ws.Range("E49").Value = "  +0.90%  "
# Row 50 - MXToken
$ws.Range("E50").Value = "  -1.01%  "
# Row 51 - RocketPoolETH
$ws.Range("E51").Value = "  +0.35%  "
